$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.747.74'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.331.60'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").Value = '3.328.97'
$ws.Range("E9").Value = '  +1.45%  '
$ws.Range("E10").Value = '  +3.77%  '
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '690.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.29%  '
$ws.Range("D15").Value = '3.872.66'
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").Value = '67.792.69'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("D19").Value = '3.325.81'
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.895'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("E30").Value = '  +1.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.09'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '571.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.99'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("E34").Value = '  +1.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.35'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = '3.712.98'
$ws.Range("E37").Value = '  -4.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("E39").Value = '  +2.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.64'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.27%  '
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("E45").Value = '  +1.97%  '
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("E47").Value = '  +5.09%  '
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("E50").Value = '  -4.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.27%  '
